$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$data = @(
    @(800, 50, 0.5, 0.6, 4648, 0),
    @(800, 50, 0.5, 0.6, 4569, -100),
    @(800, 50, 0.5, 0.6, 2773, 0),
    @(800, 50, 0.5, 0.6, 4726, 0),
    @(800, 50, 0.5, 0.6, 5464, -100),
    @(800, 50, 0.5, 0.6, 2784, 0),
    @(800, 50, 0.5, 0.6, 5126, 0),
    @(800, 50, 0.5, 0.6, 4583, 0),
    @(800, 50, 0.5, 0.6, 4590, 0),
    @(800, 50, 0.5, 0.6, 5144, 0),
    @(800, 50, 0.5, 0.6, 4397, 0),
    @(800, 50, 0.5, 0.6, 4757, 0),
    @(800, 100, 0.5, 0.6, 10327, 0),
    @(800, 100, 0.5, 0.6, 6517, -100),
    @(800, 100, 0.5, 0.7, 12071, -100),
    @(800, 70, 0.3, 0.7, 3338, 0),
    @(800, 70, 0.3, 0.7, 3236, 0),
    @(800, 70, 0.3, 0.7, 2806, 0),
    @(800, 70, 0.3, 0.7, 3941, 0),
    @(800, 70, 0.3, 0.7, 3820, 0),
    @(800, 70, 0.3, 0.7, 3734, 0),
    @(800, 70, 0.3, 0.7, 4230, 0),
    @(800, 70, 0.3, 0.7, 3706, 0),
    @(500, 70, 0.3, 0.7, 2683, 0),
    @(500, 70, 0.3, 0.7, 2699, 0),
    @(500, 70, 0.3, 0.7, 2430, 0),
    @(500, 70, 0.3, 0.7, 2765, 20),
    @(500, 70, 0.3, 0.7, 1300, 20),
    @(800, 70, 0.3, 0.7, 3609, 20),
    @(800, 70, 0.6, 0.7, 4265, 20),
    @(800, 70, 0.6, 0.7, 4130, 20),
    @(800, 70, 0.6, 0.7, 4146, 20),
    @(800, 70, 0.6, 0.7, 4342, 0),
    @(800, 70, 0.6, 0.7, 4278, 0),
    @(800, 70, 0.6, 0.7, 4156, 0),
    @(800, 70, 0.6, 0.7, 4154, 0),
    @(800, 70, 0.6, 0.7, 4293, 0),
    @(800, 70, 0.6, 0.7, 4325, 0),
    @(800, 50, 0.6, 0.7, 1289, 0),
    @(800, 50, 0.6, 0.7, 9238, 0),
    @(800, 50, 0.6, 0.7, 2603, 0),
    @(800, 50, 0.6, 0.7, 5298, 0),
    @(800, 50, 0.6, 0.7, 4145, 0),
    @(800, 50, 0.6, 0.7, 9047, 0),
    @(800, 50, 0.6, 0.7, 1366, 0),
    @(800, 50, 0.6, 0.7, 1372, 0),
    @(800, 50, 0.6, 0.7, 2134, 0),
    @(800, 50, 0.6, 0.7, 2105, 0)
)

$startRow = 253
$endRow = 300
$numRows = $endRow - $startRow + 1

$values = New-Object 'object[,]' $numRows,6
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt 6; $j++) {
        $values[$i, $j] = $data[$i][$j]
    }
}

$range = $ws.Range("A$startRow" + ":F$endRow")
$range.Value = $values
